$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.03870662259763
$ws.Range("C2").Value = 0.1005698810537581
$ws.Range("D2").Value = 0.04470433944929653
$ws.Range("E2").Value = 0.04357560434474195
$ws.Range("F2").Value = 1.356726192630845
$ws.Range("G2").Value = 1.090334264353288
$ws.Range("I2").Value = 0.03148712518962604
$ws.Range("J2").Value = 0.7731605086057556
$ws.Range("K2").Value = 0.8378806489968795
$ws.Range("L2").Value = 0.06760997928455126
$ws.Range("M2").Value = 0.9695249532526589
$ws.Range("N2").Value = 0.2102634802266437

$ws.Range("B3").Value = 0.9061340624123204
$ws.Range("C3").Value = 0.08821733117036956
$ws.Range("D3").Value = 0.0397377624895654
$ws.Range("E3").Value = 0.04177595703568926
$ws.Range("F3").Value = 1.318499581788402
$ws.Range("G3").Value = 1.059177875018349
$ws.Range("I3").Value = 0.03810478099778258
$ws.Range("J3").Value = 0.7648727961346538
$ws.Range("K3").Value = 0.8277793599572618
$ws.Range("L3").Value = 0.06727295282057355
$ws.Range("M3").Value = 0.8430345471709018
$ws.Range("N3").Value = 0.1912141442336051

$ws.Range("B4").Value = 0.8248452903971213
$ws.Range("C4").Value = 0.080765482027644
$ws.Range("D4").Value = 0.03673775110507194
$ws.Range("E4").Value = 0.04069797091330329
$ws.Range("F4").Value = 1.295969797712274
$ws.Range("G4").Value = 1.040791458748856
$ws.Range("I4").Value = 0.04267279862679807
$ws.Range("J4").Value = 0.7602580267589047
$ws.Range("K4").Value = 0.82195717096522
$ws.Range("L4").Value = 0.06708393682050051
$ws.Range("M4").Value = 0.7657295999099745
$ws.Range("N4").Value = 0.179812289410954

$ws.Range("B5").Value = 0.7914935220058794
$ws.Range("C5").Value = 0.07797099107613548
$ws.Range("D5").Value = 0.03559620595459734
$ws.Range("E5").Value = 0.04025457226414453
$ws.Range("F5").Value = 1.285704880592704
$ws.Range("G5").Value = 1.032164974233979
$ws.Range("I5").Value = 0.04476721408327178
$ws.Range("J5").Value = 0.7578195949729434
$ws.Range("K5").Value = 0.8187011216264608
$ws.Range("L5").Value = 0.06696009439120409
$ws.Range("M5").Value = 0.7345877577317594
$ws.Range("N5").Value = 0.1754051944065651

$ws.Range("B6").Value = 0.7856569416677246
$ws.Range("C6").Value = 0.07776958470805084
$ws.Range("D6").Value = 0.03549209745576576
$ws.Range("E6").Value = 0.04016825479592256
$ws.Range("F6").Value = 1.282413106473825
$ws.Range("G6").Value = 1.029141339961555
$ws.Range("I6").Value = 0.04526418983673519
$ws.Range("J6").Value = 0.7565995246469015
$ws.Range("K6").Value = 0.8169752179567098
$ws.Range("L6").Value = 0.06687720066373259
$ws.Range("M6").Value = 0.7297671780294479
$ws.Range("N6").Value = 0.1748845570654183

$ws.Range("B7").Value = 0.8235724863159817
$ws.Range("C7").Value = 0.08144144086451632
$ws.Range("D7").Value = 0.03695381663000319
$ws.Range("E7").Value = 0.04065616486229651
$ws.Range("F7").Value = 1.291461376565486
$ws.Range("G7").Value = 1.036301773813221
$ws.Range("I7").Value = 0.0430679874378832
$ws.Range("J7").Value = 0.7579829276254344
$ws.Range("K7").Value = 0.8186636900871136
$ws.Range("L7").Value = 0.06691154330802185
$ws.Range("M7").Value = 0.766254913393027
$ws.Range("N7").Value = 0.1803179329702687

$ws.Range("B8").Value = 0.9918548809792185
$ws.Range("C8").Value = 0.0972340980683839
$ws.Range("D8").Value = 0.04328847505185252
$ws.Range("E8").Value = 0.04290165998784445
$ws.Range("F8").Value = 1.337579413676835
$ws.Range("G8").Value = 1.073661225203139
$ws.Range("I8").Value = 0.034097008677624
$ws.Range("J8").Value = 0.7672533183220764
$ws.Range("K8").Value = 0.8300374218355699
$ws.Range("L8").Value = 0.06726563538128971
$ws.Range("M8").Value = 0.9270804612926327
$ws.Range("N8").Value = 0.2043807038187708

$ws.Range("B9").Value = 1.324721267603593
$ws.Range("C9").Value = 0.1281368429365557
$ws.Range("D9").Value = 0.05571627989930761
$ws.Range("E9").Value = 0.04759438797811555
$ws.Range("F9").Value = 1.443810150741641
$ws.Range("G9").Value = 1.161111321039741
$ws.Range("I9").Value = 0.02016189172708049
$ws.Range("J9").Value = 0.793383743406352
$ws.Range("K9").Value = 0.8613352545461623
$ws.Range("L9").Value = 0.06840856080814817
$ws.Range("M9").Value = 1.244253689846715
$ws.Range("N9").Value = 0.2530735433234099

$ws.Range("B10").Value = 1.570384483935328
$ws.Range("C10").Value = 0.1522109858878054
$ws.Range("D10").Value = 0.06526310595753415
$ws.Range("E10").Value = 0.05117669079665355
$ws.Range("F10").Value = 1.525399086924139
$ws.Range("G10").Value = 1.227850991991204
$ws.Range("I10").Value = 0.01294733728518604
$ws.Range("J10").Value = 0.8142325940003445
$ws.Range("K10").Value = 0.8850118755014122
$ws.Range("L10").Value = 0.06926018274080548
$ws.Range("M10").Value = 1.480824550346767
$ws.Range("N10").Value = 0.2909354927188446

$ws.Range("B11").Value = 1.681184307651165
$ws.Range("C11").Value = 0.1646623631435062
$ws.Range("D11").Value = 0.07006705352864628
$ws.Range("E11").Value = 0.05278287473955068
$ws.Range("F11").Value = 1.556888335880615
$ws.Range("G11").Value = 1.252340682967144
$ws.Range("I11").Value = 0.01083510900157236
$ws.Range("J11").Value = 0.8208143557558714
$ws.Range("K11").Value = 0.8911521449574735
$ws.Range("L11").Value = 0.06939824279035989
$ws.Range("M11").Value = 1.590885397506383
$ws.Range("N11").Value = 0.3095474547645694

$ws.Range("B12").Value = 1.724000664306431
$ws.Range("C12").Value = 0.168777643391266
$ws.Range("D12").Value = 0.07169343150499685
$ws.Range("E12").Value = 0.05342870531099031
$ws.Range("F12").Value = 1.572791491275353
$ws.Range("G12").Value = 1.265557930841112
$ws.Range("I12").Value = 0.00988311329723679
$ws.Range("J12").Value = 0.8253167870134206
$ws.Range("K12").Value = 0.8963356460502609
$ws.Range("L12").Value = 0.06959856584568769
$ws.Range("M12").Value = 1.631889482672051
$ws.Range("N12").Value = 0.3161716953413531

$ws.Range("B13").Value = 1.714931553318678
$ws.Range("C13").Value = 0.1677586015719044
$ws.Range("D13").Value = 0.07130136092037276
$ws.Range("E13").Value = 0.05329581199168842
$ws.Range("F13").Value = 1.570110778591314
$ws.Range("G13").Value = 1.263457253136252
$ws.Range("I13").Value = 0.01002974179316585
$ws.Range("J13").Value = 0.8247248406773764
$ws.Range("K13").Value = 0.895770900428623
$ws.Range("L13").Value = 0.06958429774574526
$ws.Range("M13").Value = 1.622886248797528
$ws.Range("N13").Value = 0.314642072549276

$ws.Range("B14").Value = 1.684771250992924
$ws.Range("C14").Value = 0.1649445061023442
$ws.Range("D14").Value = 0.07018308804175888
$ws.Range("E14").Value = 0.05283863213061579
$ws.Range("F14").Value = 1.558512639591711
$ws.Range("G14").Value = 1.253744713265348
$ws.Range("I14").Value = 0.01073144830426553
$ws.Range("J14").Value = 0.8213451563069896
$ws.Range("K14").Value = 0.8918128320824081
$ws.Range("L14").Value = 0.06942698769653788
$ws.Range("M14").Value = 1.594185466993252
$ws.Range("N14").Value = 0.3100486115831274

$ws.Range("B15").Value = 1.666002583881436
$ws.Range("C15").Value = 0.1634841454522444
$ws.Range("D15").Value = 0.06958093412026045
$ws.Range("E15").Value = 0.052546652881027
$ws.Range("F15").Value = 1.549952820661645
$ws.Range("G15").Value = 1.246335437924174
$ws.Range("I15").Value = 0.01128599021183163
$ws.Range("J15").Value = 0.8185355877772054
$ws.Range("K15").Value = 0.8883064466600743
$ws.Range("L15").Value = 0.06927392405026822
$ws.Range("M15").Value = 1.576951233947057
$ws.Range("N15").Value = 0.3074408434828371

$ws.Range("B16").Value = 1.560586383910646
$ws.Range("C16").Value = 0.1535595260627076
$ws.Range("D16").Value = 0.06563380832474053
$ws.Range("E16").Value = 0.05096582485116663
$ws.Range("F16").Value = 1.510833554109013
$ws.Range("G16").Value = 1.213736108750766
$ws.Range("I16").Value = 0.01399414892015738
$ws.Range("J16").Value = 0.8074522288425356
$ws.Range("K16").Value = 0.8753774881186018
$ws.Range("L16").Value = 0.06876847868199221
$ws.Range("M16").Value = 1.476398021950558
$ws.Range("N16").Value = 0.2913818817816463

$ws.Range("B17").Value = 1.496072706874429
$ws.Range("C17").Value = 0.1475366018081559
$ws.Range("D17").Value = 0.06322928761833424
$ws.Range("E17").Value = 0.0500062473554781
$ws.Range("F17").Value = 1.487268214648793
$ws.Range("G17").Value = 1.194107652395147
$ws.Range("I17").Value = 0.0158458841906528
$ws.Range("J17").Value = 0.8008570610178225
$ws.Range("K17").Value = 0.8676361279084119
$ws.Range("L17").Value = 0.06846620409006476
$ws.Range("M17").Value = 1.414943846910688
$ws.Range("N17").Value = 0.2816438383816973

$ws.Range("B18").Value = 1.459789980343061
$ws.Range("C18").Value = 0.14343997447844
$ws.Range("D18").Value = 0.06164505882208715
$ws.Range("E18").Value = 0.04949016418558116
$ws.Range("F18").Value = 1.477656592167989
$ws.Range("G18").Value = 1.186735992118116
$ws.Range("I18").Value = 0.01666522953453775
$ws.Range("J18").Value = 0.7990609539731537
$ws.Range("K18").Value = 0.8660425513362355
$ws.Range("L18").Value = 0.06844113239214167
$ws.Range("M18").Value = 1.378842892024664
$ws.Range("N18").Value = 0.2755815425005608

$ws.Range("B19").Value = 1.446918467124789
$ws.Range("C19").Value = 0.1425615006176599
$ws.Range("D19").Value = 0.06126943373968885
$ws.Range("E19").Value = 0.04929120374775664
$ws.Range("F19").Value = 1.471504428054644
$ws.Range("G19").Value = 1.181339502327205
$ws.Range("I19").Value = 0.01719233349515203
$ws.Range("J19").Value = 0.7969801384937227
$ws.Range("K19").Value = 0.8633601886232114
$ws.Range("L19").Value = 0.06832048004845603
$ws.Range("M19").Value = 1.367275408886115
$ws.Range("N19").Value = 0.2739222627994025

$ws.Range("B20").Value = 1.502995493923635
$ws.Range("C20").Value = 0.1481238923699664
$ws.Range("D20").Value = 0.06346835365832248
$ws.Range("E20").Value = 0.05011054639297008
$ws.Range("F20").Value = 1.490058032980031
$ws.Range("G20").Value = 1.196480577699162
$ws.Range("I20").Value = 0.01561701859293496
$ws.Range("J20").Value = 0.8017024645444621
$ws.Range("K20").Value = 0.868672051155734
$ws.Range("L20").Value = 0.06850953360548218
$ws.Range("M20").Value = 1.42141202572347
$ws.Range("N20").Value = 0.2826369259267523

$ws.Range("B21").Value = 1.69275172275016
$ws.Range("C21").Value = 0.1665035886788786
$ws.Range("D21").Value = 0.07074286652789397
$ws.Range("E21").Value = 0.05293657814409691
$ws.Range("F21").Value = 1.557684096730924
$ws.Range("G21").Value = 1.252364454825241
$ws.Range("I21").Value = 0.0108356606539548
$ws.Range("J21").Value = 0.8201913092784849
$ws.Range("K21").Value = 0.8898598403244549
$ws.Range("L21").Value = 0.06931047271246271
$ws.Range("M21").Value = 1.603539302730695
$ws.Range("N21").Value = 0.3119554604457591

$ws.Range("B22").Value = 1.818452589071683
$ws.Range("C22").Value = 0.1778181329958386
$ws.Range("D22").Value = 0.07526133189590212
$ws.Range("E22").Value = 0.05486430312029356
$ws.Range("F22").Value = 1.608745277973469
$ws.Range("G22").Value = 1.295538438045511
$ws.Range("I22").Value = 0.007956544520742348
$ws.Range("J22").Value = 0.8356960154757047
$ws.Range("K22").Value = 0.9083134635349026
$ws.Range("L22").Value = 0.07006669597514303
$ws.Range("M22").Value = 1.722228667139262
$ws.Range("N22").Value = 0.3307993941624545

$ws.Range("B23").Value = 1.752364455797533
$ws.Range("C23").Value = 0.170872738413081
$ws.Range("D23").Value = 0.07256512556497796
$ws.Range("E23").Value = 0.05387579502938067
$ws.Range("F23").Value = 1.58645027613845
$ws.Range("G23").Value = 1.277472427132906
$ws.Range("I23").Value = 0.009073446207194635
$ws.Range("J23").Value = 0.8299391804437306
$ws.Range("K23").Value = 0.9021578156062375
$ws.Range("L23").Value = 0.06985678786664096
$ws.Range("M23").Value = 1.657681204638749
$ws.Range("N23").Value = 0.3200300822091151

$ws.Range("B24").Value = 1.501503419189987
$ws.Range("C24").Value = 0.1464749658996851
$ws.Range("D24").Value = 0.06292221472094894
$ws.Range("E24").Value = 0.05013160712837816
$ws.Range("F24").Value = 1.496824562675741
$ws.Range("G24").Value = 1.203434473812493
$ws.Range("I24").Value = 0.01510045091966106
$ws.Range("J24").Value = 0.8053970322495587
$ws.Range("K24").Value = 0.8741212448263767
$ws.Range("L24").Value = 0.06879930428408443
$ws.Range("M24").Value = 1.416725450522961
$ws.Range("N24").Value = 0.2811277240044205

$ws.Range("B25").Value = 1.23289742883324
$ws.Range("C25").Value = 0.1209437842112067
$ws.Range("D25").Value = 0.05272747088713459
$ws.Range("E25").Value = 0.04623349867360638
$ws.Range("F25").Value = 1.406298760831106
$ws.Range("G25").Value = 1.128855481143006
$ws.Range("I25").Value = 0.02402661446872045
$ws.Range("J25").Value = 0.7818654015633228
$ws.Range("K25").Value = 0.8466971654974316
$ws.Range("L25").Value = 0.06777980217629676
$ws.Range("M25").Value = 1.159694507752562
$ws.Range("N25").Value = 0.2406307789391349
